# Apply changes described by the commit "fix tariffer and setting docstirngs"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet (timestamp updated from 15_02 to 16_44)
$ws.Name = "2023_07_05 16_44"

# Row 2
$ws.Range("D2").Value2 = -64
$ws.Range("I2").Value2 = 24597
$ws.Range("J2").Value2 = 40347
$ws.Range("K2").Value2 = 24509

# Row 3
$ws.Range("D3").Value2 = -58
$ws.Range("I3").Value2 = 23773
$ws.Range("J3").Value2 = 38686
$ws.Range("K3").Value2 = 23655

# Row 4
$ws.Range("D4").Value2 = -64
$ws.Range("G4").Value2 = 26
$ws.Range("I4").Value2 = 24526
$ws.Range("J4").Value2 = 40215
$ws.Range("K4").Value2 = 24420

# Row 5
$ws.Range("D5").Value2 = -74
$ws.Range("I5").Value2 = 24571
$ws.Range("J5").Value2 = 40303
$ws.Range("K5").Value2 = 24471

# Row 6
$ws.Range("D6").Value2 = -66
$ws.Range("G6").Value2 = 25
$ws.Range("I6").Value2 = 24471
$ws.Range("J6").Value2 = 40099

# Row 7
$ws.Range("D7").Value2 = -62
$ws.Range("G7").Value2 = 23
$ws.Range("I7").Value2 = 24613
$ws.Range("J7").Value2 = 40380
$ws.Range("K7").Value2 = 24525

# Row 8
$ws.Range("D8").Value2 = -60
$ws.Range("G8").Value2 = 24
$ws.Range("I8").Value2 = 24542
$ws.Range("J8").Value2 = 40239
$ws.Range("K8").Value2 = 24448

# Row 9
$ws.Range("D9").Value2 = -66
$ws.Range("I9").Value2 = 24205
$ws.Range("J9").Value2 = 39552
$ws.Range("K9").Value2 = 24111

# Row 10
$ws.Range("G10").Value2 = 26
$ws.Range("I10").Value2 = 23865
$ws.Range("J10").Value2 = 38866
$ws.Range("K10").Value2 = 23759

# Row 11
$ws.Range("D11").Value2 = -60
$ws.Range("I11").Value2 = 24311
$ws.Range("J11").Value2 = 39772
$ws.Range("K11").Value2 = 24211
